$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the dataset. It belongs right after the
# header block that starts the "Nectarín" rows tracked on this sheet, at
# row 21, so insert a fresh row there (this shifts the former rows 21-34
# down to 22-35, carrying their data and formatting with them) and then
# populate the new row 21 with the new record's values.
$ws.Rows(21).Insert()

$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = 44523
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100103
$ws.Cells.Item(21, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(21, 9).Value = 100103006
$ws.Cells.Item(21, 10).Value = "Nectarín"
$ws.Cells.Item(21, 11).Value = "Early Glo"
$ws.Cells.Item(21, 12).Value = "Segunda"
$ws.Cells.Item(21, 13).Value = 250
$ws.Cells.Item(21, 14).Value = 24000
$ws.Cells.Item(21, 15).Value = 25000
$ws.Cells.Item(21, 16).Value = 24500
$ws.Cells.Item(21, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(21, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(21, 19).Value = 1361
$ws.Cells.Item(21, 20).Value = 18
